$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Repull / push new dSF (column F) values for rows 7-35 (mean calculation refresh)
$updates = @{
    7  = -6
    8  = 4
    9  = -7
    10 = -2
    11 = -4
    12 = 2
    13 = 7
    14 = 2
    15 = -3
    16 = 3
    17 = 6
    18 = 6
    19 = -2
    20 = 1
    21 = -1
    22 = -1
    23 = -2
    25 = 2
    26 = 0
    27 = 4
    28 = -1
    29 = 5
    30 = -3
    32 = -5
    33 = -2
    35 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
